$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "NAME"
$ws.Range("B1").Value = "Dec/Hep"
$ws.Range("C1").Value = "Result"
$ws.Range("D1").Value = "Score"
$ws.Range("E1").Value = "Total"

# Row 2 - sten
$ws.Range("A2").Value = "sten"
$ws.Range("B2").Value = "Hep 200M"
$ws.Range("C2").Value = 22.0
$ws.Range("D2").Value = 1181.0
$ws.Range("E2").Value = 1181.0

# Row 3 - Flisa
$ws.Range("A3").Value = "Flisa"
$ws.Range("B3").Value = "Hep 200M"
$ws.Range("C3").Value = 25.0
$ws.Range("D3").Value = 887.0
$ws.Range("E3").Value = 887.0

# Row 4 - sten
$ws.Range("A4").Value = "sten"
$ws.Range("B4").Value = "Hep 800M"
$ws.Range("C4").Value = 98.0
$ws.Range("D4").Value = 1486.0
$ws.Range("E4").Value = 2667.0
